$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the music value cells (B24, B33, B39, B44, B50) from "_" to "bgm_battle_a1"
$musicRows = @(24, 33, 39, 44, 50)
foreach ($r in $musicRows) {
    $ws.Range("B$r").Value = "bgm_battle_a1"
}

# Update the sheet view: scroll the visible window so row 40 is the
# top-left row (topLeftCell A40) and select B50 (matching the new
# activeCell/sqref in the diff).
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("B50").Select()
